$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# --- Step 1: Update Ybus-derived mode values (hybrid ac/dc network) ---
$ws1 = $wb.Worksheets.Item("State-PF")
$ws1.Range("F11").Value = "-2494320.45+50.00i"
$ws1.Range("F12").Value = "-2494320.45-50.00i"
$ws1.Range("F13").Value = "-230.40+135868.84i"
$ws1.Range("F14").Value = "-230.40-135868.84i"
$ws1.Range("F15").Value = "-230.40+135768.84i"
$ws1.Range("F16").Value = "-230.40-135768.84i"
$ws1.Range("F17").Value = "-18.65+119506.08i"
$ws1.Range("F18").Value = "-18.65-119506.08i"
$ws1.Range("F19").Value = "-18.65+119406.08i"
$ws1.Range("F20").Value = "-18.65-119406.08i"
$ws1.Range("F21").Value = "-29.14+98918.21i"
$ws1.Range("F22").Value = "-29.14-98918.21i"
$ws1.Range("F23").Value = "-29.14+99018.21i"
$ws1.Range("F24").Value = "-29.14-99018.21i"
$ws1.Range("F25").Value = "-5661.60+50.00i"
$ws1.Range("F26").Value = "-5661.60-50.00i"
$ws1.Range("F28").Value = "-25.78+100.29i"
$ws1.Range("F29").Value = "-25.78-100.29i"
$ws1.Range("F30").Value = "-3.14+54.10i"
$ws1.Range("F31").Value = "-3.14-54.10i"
$ws1.Range("F32").Value = "-5.56+50.00i"
$ws1.Range("F33").Value = "-5.56-50.00i"
$ws1.Range("F34").Value = "-5.65+50.04i"
$ws1.Range("F35").Value = "-5.65-50.04i"
$ws1.Range("F37").Value = "-0.08+8.13i"
$ws1.Range("F38").Value = "-0.08-8.13i"
$ws1.Range("F39").Value = "-6.99+4.77i"
$ws1.Range("F40").Value = "-6.99-4.77i"
$ws1.Range("F42").Value = "-0.00+3.01i"
$ws1.Range("F43").Value = "-0.00-3.01i"
$ws1.Range("F36").Value = -17.9
$ws1.Range("F41").Value = -6.52
$ws1.Range("F45").Value = -0.04

$ws2 = $wb.Worksheets.Item("Impedance-PF")
$ws2.Range("H12").Value = "-2494320.45+50.00i"
$ws2.Range("H13").Value = "-2494320.45-50.00i"
$ws2.Range("H14").Value = "-230.40+135868.84i"
$ws2.Range("H15").Value = "-230.40-135868.84i"
$ws2.Range("H16").Value = "-230.40+135768.84i"
$ws2.Range("H17").Value = "-230.40-135768.84i"
$ws2.Range("H18").Value = "-18.65+119506.08i"
$ws2.Range("H19").Value = "-18.65-119506.08i"
$ws2.Range("H20").Value = "-18.65+119406.08i"
$ws2.Range("H21").Value = "-18.65-119406.08i"
$ws2.Range("H22").Value = "-29.14+98918.21i"
$ws2.Range("H23").Value = "-29.14-98918.21i"
$ws2.Range("H24").Value = "-29.14+99018.21i"
$ws2.Range("H25").Value = "-29.14-99018.21i"
$ws2.Range("H26").Value = "-5661.60+50.00i"
$ws2.Range("H27").Value = "-5661.60-50.00i"
$ws2.Range("H29").Value = "-25.78+100.29i"
$ws2.Range("H30").Value = "-25.78-100.29i"
$ws2.Range("H31").Value = "-3.14+54.10i"
$ws2.Range("H32").Value = "-3.14-54.10i"
$ws2.Range("H33").Value = "-5.56+50.00i"
$ws2.Range("H34").Value = "-5.56-50.00i"
$ws2.Range("H35").Value = "-5.65+50.04i"
$ws2.Range("H36").Value = "-5.65-50.04i"
$ws2.Range("H38").Value = "-0.08+8.13i"
$ws2.Range("H39").Value = "-0.08-8.13i"
$ws2.Range("H40").Value = "-6.99+4.77i"
$ws2.Range("H41").Value = "-6.99-4.77i"
$ws2.Range("H43").Value = "-0.00+3.01i"
$ws2.Range("H44").Value = "-0.00-3.01i"
$ws2.Range("H37").Value = -17.9
$ws2.Range("H42").Value = -6.52
$ws2.Range("H46").Value = -0.04

# --- Step 2: Bump internal sheetId bookkeeping (State-PF/Impedance-PF/Enabling
#     move from sheetId 123/124/125 to 141/142/143). Excel assigns each new
#     sheet the next sheetId = (max sheetId currently in the workbook) + 1, so
#     we temporarily pad the workbook with throwaway sheets to raise that
#     ceiling, duplicate the three real sheets (which then inherit 141/142/143
#     in sequence), and drop the scaffolding again. ---
for ($i = 1; $i -le 15; $i++) {
    $dummy = $wb.Worksheets.Add()
    $dummy.Name = "__tmp$i"
}

$lastIndex = $wb.Worksheets.Count
$wb.Worksheets.Item("State-PF").Copy($null, $wb.Worksheets.Item($lastIndex))
$lastIndex = $wb.Worksheets.Count
$wb.Worksheets.Item("Impedance-PF").Copy($null, $wb.Worksheets.Item($lastIndex))
$lastIndex = $wb.Worksheets.Count
$wb.Worksheets.Item("Enabling").Copy($null, $wb.Worksheets.Item($lastIndex))

for ($i = 1; $i -le 15; $i++) {
    $wb.Worksheets.Item("__tmp$i").Delete()
}

$wb.Worksheets.Item("State-PF").Delete()
$wb.Worksheets.Item("Impedance-PF").Delete()
$wb.Worksheets.Item("Enabling").Delete()

$wb.Worksheets.Item("State-PF (2)").Name = "State-PF"
$wb.Worksheets.Item("Impedance-PF (2)").Name = "Impedance-PF"
$wb.Worksheets.Item("Enabling (2)").Name = "Enabling"

# Restore original tab order / active tab ("Enabling" selected, as before)
$wb.Worksheets.Item("State-PF").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("Impedance-PF").Move($null, $wb.Worksheets.Item("State-PF"))
$wb.Worksheets.Item("Enabling").Move($null, $wb.Worksheets.Item("Impedance-PF"))
$wb.Worksheets.Item("Enabling").Activate()
